# Apply edit: rotate three swear-word cells and update the active selection
# to show the "window" (view) scrolled/selected at E23, per commit message
# "added window to show censored words".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rotate values of A19, A20, A21:
#   A19 <- old A20
#   A20 <- old A21
#   A21 <- old A19
$v19 = $ws.Range("A19").Value2
$v20 = $ws.Range("A20").Value2
$v21 = $ws.Range("A21").Value2

$ws.Range("A19").Value2 = $v20
$ws.Range("A20").Value2 = $v21
$ws.Range("A21").Value2 = $v19

# Update the active cell / selection shown in the sheet view.
$ws.Range("E23").Select()
